$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3256646321388814
$ws.Range("D2").Value = 0.01125752210194619
$ws.Range("E2").Value = 0.1217552071829093
$ws.Range("F2").Value = 0.7931065201222935
$ws.Range("G2").Value = 0.6514338298070754
$ws.Range("H2").Value = 0.6938152394459394
$ws.Range("L2").Value = 0.1047364062361815
$ws.Range("M2").Value = 1.098681838835375
$ws.Range("N2").Value = 1.793565731190284
$ws.Range("O2").Value = 2.679343611211209

$ws.Range("C3").Value = 0.3269899379945969
$ws.Range("D3").Value = 0.0111022946356627
$ws.Range("E3").Value = 0.1240207849380939
$ws.Range("F3").Value = 0.7703493223745141
$ws.Range("G3").Value = 0.6277865032679273
$ws.Range("H3").Value = 0.6875998781628567
$ws.Range("L3").Value = 0.1066407583515812
$ws.Range("M3").Value = 0.9873498382050201
$ws.Range("N3").Value = 1.648381501539035
$ws.Range("O3").Value = 2.615801535586371

$ws.Range("C4").Value = 0.3280252179163341
$ws.Range("D4").Value = 0.01100670368832724
$ws.Range("E4").Value = 0.1255002352143997
$ws.Range("F4").Value = 0.7569334786185493
$ws.Range("G4").Value = 0.6137601379654853
$ws.Range("H4").Value = 0.6841965732360364
$ws.Range("L4").Value = 0.1078781147655175
$ws.Range("M4").Value = 0.9188200035312377
$ws.Range("N4").Value = 1.559417443217882
$ws.Range("O4").Value = 2.578681299817674

$ws.Range("C5").Value = 0.3285026734227614
$ws.Range("D5").Value = 0.01096768421183825
$ws.Range("E5").Value = 0.1261252826731125
$ws.Range("F5").Value = 0.751606110585783
$ws.Range("G5").Value = 0.6081677544764972
$ws.Range("H5").Value = 0.6829134068727996
$ws.Range("L5").Value = 0.1083994460922746
$ws.Range("M5").Value = 0.8908530145642288
$ws.Range("N5").Value = 1.523213347888429
$ws.Range("O5").Value = 2.564029542699529

$ws.Range("C6").Value = 0.3285853068774856
$ws.Range("D6").Value = 0.010961201261658
$ws.Range("E6").Value = 0.1262304078454166
$ws.Range("F6").Value = 0.7507299307017377
$ws.Range("G6").Value = 0.6072465851927262
$ws.Range("H6").Value = 0.6827065979625786
$ws.Range("L6").Value = 0.1084870448950515
$ws.Range("M6").Value = 0.8862067467534729
$ws.Range("N6").Value = 1.517204801562428
$ws.Range("O6").Value = 2.561625279032029

$ws.Range("C7").Value = 0.3280314322257993
$ws.Range("D7").Value = 0.01100617771652956
$ws.Range("E7").Value = 0.1255085751847158
$ws.Range("F7").Value = 0.7568610667830455
$ws.Range("G7").Value = 0.6136842179023319
$ws.Range("H7").Value = 0.6841788482743283
$ws.Range("L7").Value = 0.1078850764178281
$ws.Range("M7").Value = 0.9184429915554801
$ws.Range("N7").Value = 1.558928976037038
$ws.Range("O7").Value = 2.578481779589254

$ws.Range("C8").Value = 0.3260755347638025
$ws.Range("D8").Value = 0.01120406096328885
$ws.Range("E8").Value = 0.1225179979496591
$ws.Range("F8").Value = 0.7851439437247762
$ws.Range("G8").Value = 0.6431775349374789
$ws.Range("H8").Value = 0.6915863909422484
$ws.Range("L8").Value = 0.1053788857672817
$ws.Range("M8").Value = 1.060331694352371
$ws.Range("N8").Value = 1.743471444944987
$ws.Range("O8").Value = 2.657040133589987

$ws.Range("C9").Value = 0.3240042782072408
$ws.Range("D9").Value = 0.01158967540727573
$ws.Range("E9").Value = 0.1173574690164916
$ws.Range("F9").Value = 0.8450497950826872
$ws.Range("G9").Value = 0.7049553397110628
$ws.Range("H9").Value = 0.709396879391619
$ws.Range("L9").Value = 0.1010052377591215
$ws.Range("M9").Value = 1.337112848595339
$ws.Range("N9").Value = 2.106614947438516
$ws.Range("O9").Value = 2.826204250012438

$ws.Range("C10").Value = 0.3235674214380992
$ws.Range("D10").Value = 0.01187125436973346
$ws.Range("E10").Value = 0.1139988574648036
$ws.Range("F10").Value = 0.8918079835397492
$ws.Range("G10").Value = 0.752790963821667
$ws.Range("H10").Value = 0.7244977131446433
$ws.Range("L10").Value = 0.09812276530411523
$ws.Range("M10").Value = 1.539456627227068
$ws.Range("N10").Value = 2.37398452997229
$ws.Range("O10").Value = 2.959824487267838

$ws.Range("C11").Value = 0.3236063077427076
$ws.Range("D11").Value = 0.01199892609499287
$ws.Range("E11").Value = 0.1125657059727239
$ws.Range("F11").Value = 0.9136838764056563
$ws.Range("G11").Value = 0.7750940202518564
$ws.Range("H11").Value = 0.7318081516062023
$ws.Range("L11").Value = 0.09688352189750304
$ws.Range("M11").Value = 1.631266146435351
$ws.Range("N11").Value = 2.495701865710714
$ws.Range("O11").Value = 3.022666468516093

$ws.Range("C12").Value = 0.3236553619494487
$ws.Range("D12").Value = 0.01204720702078887
$ws.Range("E12").Value = 0.1120366948078901
$ws.Range("F12").Value = 0.9220553050908507
$ws.Range("G12").Value = 0.783618290089521
$ws.Range("H12").Value = 0.734640036999366
$ws.Range("L12").Value = 0.09642463015339331
$ws.Range("M12").Value = 1.665995493411202
$ws.Range("N12").Value = 2.54180213172998
$ws.Range("O12").Value = 3.046760836045053

$ws.Range("C13").Value = 0.3236432679117058
$ws.Range("D13").Value = 0.01203681185804584
$ws.Range("E13").Value = 0.1121500166705092
$ws.Range("F13").Value = 0.9202484691975457
$ws.Range("G13").Value = 0.7817789314794084
$ws.Range("H13").Value = 0.7340273093564633
$ws.Range("L13").Value = 0.09652299846802137
$ws.Range("M13").Value = 1.658517592968195
$ws.Range("N13").Value = 2.531873302547183
$ws.Range("O13").Value = 3.041558426959682

$ws.Range("C14").Value = 0.3236096545435458
$ws.Range("D14").Value = 0.01200289953896316
$ws.Range("E14").Value = 0.1125219090647676
$ws.Range("F14").Value = 0.9143708425951615
$ws.Range("G14").Value = 0.7757937386817275
$ws.Range("H14").Value = 0.7320398573035618
$ws.Range("L14").Value = 0.09684556034711278
$ws.Range("M14").Value = 1.634124106424579
$ws.Range("N14").Value = 2.49949441552809
$ws.Range("O14").Value = 3.024642755575428

$ws.Range("C15").Value = 0.3235935407664243
$ws.Range("D15").Value = 0.01198211857109399
$ws.Range("E15").Value = 0.1127514892455563
$ws.Range("F15").Value = 0.9107820355081486
$ws.Range("G15").Value = 0.7721378871143827
$ws.Range("H15").Value = 0.7308307706301491
$ws.Range("L15").Value = 0.09704449180687647
$ws.Range("M15").Value = 1.619177501565929
$ws.Range("N15").Value = 2.479662398605569
$ws.Range("O15").Value = 3.014320206145214

$ws.Range("C16").Value = 0.323569675219602
$ws.Range("D16").Value = 0.01186290182552341
$ws.Range("E16").Value = 0.1140944292346324
$ws.Range("F16").Value = 0.8903905625460027
$ws.Range("G16").Value = 0.7513443714357493
$ws.Range("H16").Value = 0.7240288469284621
$ws.Range("L16").Value = 0.09820520415882328
$ws.Range("M16").Value = 1.533451641228893
$ws.Range("N16").Value = 2.366031447033038
$ws.Range("O16").Value = 2.955759155130806

$ws.Range("C17").Value = 0.3236160070926104
$ws.Range("D17").Value = 0.011789654931043
$ws.Range("E17").Value = 0.1149425869400522
$ws.Range("F17").Value = 0.8780364230754003
$ws.Range("G17").Value = 0.7387275626851419
$ws.Range("H17").Value = 0.7199691541858897
$ws.Range("L17").Value = 0.09893572647128579
$ws.Range("M17").Value = 1.480798766354724
$ws.Range("N17").Value = 2.296342383404124
$ws.Range("O17").Value = 2.920361958488229

$ws.Range("C18").Value = 0.3236650117817561
$ws.Range("D18").Value = 0.01174748603218134
$ws.Range("E18").Value = 0.1154393337368829
$ws.Range("F18").Value = 0.870987600827732
$ws.Range("G18").Value = 0.7315217249786485
$ws.Range("H18").Value = 0.7176756246110472
$ws.Range("L18").Value = 0.09936267958444311
$ws.Range("M18").Value = 1.450492016777929
$ws.Range("N18").Value = 2.256267763270557
$ws.Range("O18").Value = 2.900195988774328

$ws.Range("C19").Value = 0.3236854386886279
$ws.Range("D19").Value = 0.0117332017902676
$ws.Range("E19").Value = 0.1156090515069712
$ws.Range("F19").Value = 0.8686107585181446
$ws.Range("G19").Value = 0.7290906963713439
$ws.Range("H19").Value = 0.7169061972325892
$ws.Range("L19").Value = 0.09950840150415807
$ws.Range("M19").Value = 1.440226941783024
$ws.Range("N19").Value = 2.242700810086319
$ws.Range("O19").Value = 2.893401333308702

$ws.Range("C20").Value = 0.3236087600837578
$ws.Range("D20").Value = 0.01179745628564532
$ws.Range("E20").Value = 0.1148513765232364
$ws.Range("F20").Value = 0.8793456447120036
$ws.Range("G20").Value = 0.7400653591692503
$ws.Range("H20").Value = 0.7203970192811084
$ws.Range("L20").Value = 0.09885725962883996
$ws.Range("M20").Value = 1.48640607322973
$ws.Range("N20").Value = 2.303760035741561
$ws.Range("O20").Value = 2.924110005931993

$ws.Range("C21").Value = 0.3236185946527854
$ws.Range("D21").Value = 0.01201286223054154
$ws.Range("E21").Value = 0.1124123030544153
$ws.Range("F21").Value = 0.9160948655848102
$ws.Range("G21").Value = 0.7775495973579609
$ws.Range("H21").Value = 0.732621893241884
$ws.Range("L21").Value = 0.09675053409068468
$ws.Range("M21").Value = 1.641290088677991
$ws.Range("N21").Value = 2.50900468162348
$ws.Range("O21").Value = 3.029603212090478

$ws.Range("C22").Value = 0.3238251912401893
$ws.Range("D22").Value = 0.01215325782645138
$ws.Range("E22").Value = 0.11089807099982
$ws.Range("F22").Value = 0.9406228505335434
$ws.Range("G22").Value = 0.8025060510166213
$ws.Range("H22").Value = 0.7409822322664468
$ws.Range("L22").Value = 0.09543419963078037
$ws.Range("M22").Value = 1.742299254765271
$ws.Range("N22").Value = 2.643191803965806
$ws.Range("O22").Value = 3.10028377393337

$ws.Range("C23").Value = 0.3236965582623412
$ws.Range("D23").Value = 0.01207836294711839
$ws.Range("E23").Value = 0.1116989147117771
$ws.Range("F23").Value = 0.9274849658694109
$ws.Range("G23").Value = 0.7891441884279971
$ws.Range("H23").Value = 0.7364861895381125
$ws.Range("L23").Value = 0.09613120397173525
$ws.Range("M23").Value = 1.688409473802551
$ws.Range("N23").Value = 2.571570645921668
$ws.Range("O23").Value = 3.062400972756734

$ws.Range("C24").Value = 0.323611966794644
$ws.Range("D24").Value = 0.01179392947449287
$ws.Range("E24").Value = 0.1148925843193922
$ws.Range("F24").Value = 0.8787535782530824
$ws.Range("G24").Value = 0.739460392735225
$ws.Range("H24").Value = 0.7202034555027126
$ws.Range("L24").Value = 0.09889271278996681
$ws.Range("M24").Value = 1.483871121437147
$ws.Range("N24").Value = 2.300406544217992
$ws.Range("O24").Value = 2.922414939548673

$ws.Range("C25").Value = 0.3243747667145414
$ws.Range("D25").Value = 0.01148564298446964
$ws.Range("E25").Value = 0.1186777847563238
$ws.Range("F25").Value = 0.8283642148182651
$ws.Range("G25").Value = 0.6878162344786176
$ws.Range("H25").Value = 0.70422569664305
$ws.Range("L25").Value = 0.1021304023746286
$ws.Range("M25").Value = 1.262404802415418
$ws.Range("N25").Value = 2.008261178884936
$ws.Range("O25").Value = 2.778810398481596
